$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Split the single "SYAHPC" sheet into three sheets: green / low
#    carbon / unspecified. The original (non-zero) production data
#    is preserved on the new "unspecified" sheet; the renamed
#    "green" sheet and the new "lowcarbon" sheet start from zero.
# ------------------------------------------------------------------
$wsGreen = $wb.Worksheets.Item("SYAHPC")
$wsGreen.Name = "SYAHPC-green"

$wsLow = $wb.Worksheets.Add($null, $wsGreen)
$wsLow.Name = "SYAHPC-lowcarbon"

$wsUnspec = $wb.Worksheets.Add($null, $wsLow)
$wsUnspec.Name = "SYAHPC-unspecified"

# Copy the original layout/format/values (incl. the non-zero
# electrolysis / natural-gas-reforming capacities) onto the two new
# sheets before the "green" sheet's numbers get zeroed out below.
$wsGreen.Range("A1:B8").Copy($wsLow.Range("A1"))
$wsGreen.Range("A1:B8").Copy($wsUnspec.Range("A1"))

# Match the tab color and column width used on the original/green sheet.
$wsLow.Tab.Color = $wsGreen.Tab.Color
$wsUnspec.Tab.Color = $wsGreen.Tab.Color

$wsLow.Columns("A").ColumnWidth = $wsGreen.Columns("A").ColumnWidth
$wsUnspec.Columns("A").ColumnWidth = $wsGreen.Columns("A").ColumnWidth

[void]$wsLow.Range("B4").Select()
[void]$wsUnspec.Range("B4").Select()

# Zero out the production-capacity numbers for green + low-carbon
# (the "unspecified" sheet keeps the original historical values).
$wsGreen.Range("B2").Value = 0
$wsGreen.Range("B3").Value = 0
$wsLow.Range("B2").Value = 0
$wsLow.Range("B3").Value = 0

# ------------------------------------------------------------------
# 2. Update the "About" sheet: turn the single title line into three
#    (one per new sheet) and push the rest of the content down.
# ------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

# Old row 1 becomes rows 1-3 (titles); everything from the old row 3
# onward shifts down by two rows, so insert two rows right after the
# title row.
$wsAbout.Rows("2:3").Insert()

$wsAbout.Range("A1").Value = "SYAHPC Start Year Annual Green Hydrogen Production Capacity"
$wsAbout.Range("A2").Value = "SYAHPC Start Year Annual Low Carbon Hydrogen Production Capacity"
$wsAbout.Range("A3").Value = "SYAHPC Start Year Annual Unspecified Hydrogen Production Capacity"

[void]$wsAbout.Range("B24").Select()

# Keep "About" the active/selected sheet, as it was originally.
[void]$wsAbout.Activate()

Write-Output "done"
